# Auto-generated COM-interop script: adds "Sheet2" after "Sheet1" to Plan.xlsx,
# populating its content/styles per the target diff ("Added page to plan").

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Insert the new sheet directly after Sheet1 (mirrors Excel's "Insert Sheet" at that tab).
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Sheet2"

$v1 = @'
'- charting: can etoro take an overlay? Eg we produce a txt file, upload it and it produces a series of points on the chart that are the prediction points
'@
$ws2.Range("A1").Value = $v1

$v2 = @'
'- if not or it’s deemed too impractical, need to look into developing our own live charting or regularly updated charts
'@
$ws2.Range("A2").Value = $v2

$v3 = @'
'- our own time and availability to trade: we both have full time jobs, lives etc, how will this factor in - we need to produce an output that we can use to trade given our time availabilities
'@
$ws2.Range("A3").Value = $v3

$v4 = @'
'- keep output simple and decision making complex Eg model simply advises with its predictions we then do the decision making
'@
$ws2.Range("A4").Value = $v4

$v6 = @'
BUY/SELL INDICATORS
'@
$ws2.Range("A6").Value = $v6

$v7 = @'
Accuracy so far this month
'@
$ws2.Range("A7").Value = $v7

$v8 = @'
Cumulative value accuracy
'@
$ws2.Range("A8").Value = $v8

$v9 = @'
Averaged value accuract
'@
$ws2.Range("A9").Value = $v9

$v10 = @'
Cumulative delta accuracy
'@
$ws2.Range("A10").Value = $v10

$v11 = @'
Averaged delta accuracy
'@
$ws2.Range("A11").Value = $v11

$v12 = @'
Averaged st. dev
'@
$ws2.Range("A12").Value = $v12

$v13 = @'
no. points in boxplot
'@
$ws2.Range("A13").Value = $v13

$v14 = @'
no. points in box
'@
$ws2.Range("A14").Value = $v14

$v15 = @'
no. correct directions
'@
$ws2.Range("A15").Value = $v15

$v16 = @'
Local Accuracy
'@
$ws2.Range("A16").Value = $v16

$v17 = @'
Montly accuracy indicators for local dataset, with focus on delta values (not absolute error)
'@
$ws2.Range("A17").Value = $v17

$v18 = @'
Averaged value accuracy gradient
'@
$ws2.Range("A18").Value = $v18

$v19 = @'
Absolute value gradient
'@
$ws2.Range("A19").Value = $v19

$v20 = @'
no. correct directions
'@
$ws2.Range("A20").Value = $v20

$v21 = @'
Buy/Sell deltaT
'@
$ws2.Range("A21").Value = $v21

$v22 = @'
Potential gain/loss
'@
$ws2.Range("A22").Value = $v22

$v23 = @'
Overnight?
'@
$ws2.Range("A23").Value = $v23

$v24 = @'
Delta t length
'@
$ws2.Range("A24").Value = $v24

# Bold the section-header rows (reuses the existing bold style, s="1").
$ws2.Range("A6").Font.Bold = $true
$ws2.Range("A7").Font.Bold = $true
$ws2.Range("A16").Font.Bold = $true
$ws2.Range("A21").Font.Bold = $true

# Page setup to match: portrait orientation, default margins (already default).
$ws2.PageSetup.Orientation = 1

# Select F12 (matches the saved selection/active cell on the new sheet)
# and make Sheet2 the active/visible tab (matches activeTab + tabSelected move).
$ws2.Range("F12").Select()
$ws2.Activate()

